$wb = $excel.ActiveWorkbook

# --- Sheet "ATS Accuracy" ---
$ws1 = $wb.Worksheets.Item("ATS Accuracy")

$ws1.Range("B2").Value = 2
$ws1.Range("C2").Value = 69
$ws1.Range("D2").Value = 71
$ws1.Range("E2").Value = 97.2

$ws1.Range("B3").Value = 3
$ws1.Range("C3").Value = 57
$ws1.Range("D3").Value = 60
$ws1.Range("E3").Value = 95

$ws1.Range("B4").Value = 3
$ws1.Range("C4").Value = 13
$ws1.Range("D4").Value = 16
$ws1.Range("E4").Value = 81.2

$ws1.Range("B5").Value = 3
$ws1.Range("C5").Value = 7
$ws1.Range("D5").Value = 10
$ws1.Range("E5").Value = 70

$ws1.Range("B6").Value = 4
$ws1.Range("C6").Value = 1
$ws1.Range("D6").Value = 5
$ws1.Range("E6").Value = 20

# --- Sheet "Total Accuracy" ---
$ws2 = $wb.Worksheets.Item("Total Accuracy")

$ws2.Range("C2").Value = 62
$ws2.Range("D2").Value = 66
$ws2.Range("E2").Value = 93.90000000000001

$ws2.Range("C3").Value = 55
$ws2.Range("D3").Value = 58
$ws2.Range("E3").Value = 94.8

$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = 22
$ws2.Range("D4").Value = 24
$ws2.Range("E4").Value = 91.7

$ws2.Range("C5").Value = 8
$ws2.Range("D5").Value = 11
$ws2.Range("E5").Value = 72.7

$ws2.Range("B6").Value = 2
$ws2.Range("D6").Value = 3
$ws2.Range("E6").Value = 33.3
